$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (17 players -> 16 players; Kentavious Caldwell-Pope removed)
$ws.Rows.Item(18).Delete()

# Rewrite the full player table (rows reordered + position updates for Monk/Miller/Randle)
$ws.Cells.Item(2, 1).Value = "Cade Cunningham"
$ws.Cells.Item(2, 2).Value = "PG,SG"
$ws.Cells.Item(2, 3).Value = "Detroit Pistons"
$ws.Cells.Item(3, 1).Value = "Malik Monk"
$ws.Cells.Item(3, 2).Value = "PG,SG,SF"
$ws.Cells.Item(3, 3).Value = "Sacramento Kings"
$ws.Cells.Item(4, 1).Value = "Brandon Miller"
$ws.Cells.Item(4, 2).Value = "SG,SF,PF"
$ws.Cells.Item(4, 3).Value = "Charlotte Hornets"
$ws.Cells.Item(5, 1).Value = "Duncan Robinson"
$ws.Cells.Item(5, 2).Value = "SG,SF"
$ws.Cells.Item(5, 3).Value = "Miami Heat"
$ws.Cells.Item(6, 1).Value = "Cameron Johnson"
$ws.Cells.Item(6, 2).Value = "SF,PF"
$ws.Cells.Item(6, 3).Value = "Brooklyn Nets"
$ws.Cells.Item(7, 1).Value = "Kelly Olynyk"
$ws.Cells.Item(7, 2).Value = "C"
$ws.Cells.Item(7, 3).Value = "Toronto Raptors"
$ws.Cells.Item(8, 1).Value = "Bam Adebayo"
$ws.Cells.Item(8, 2).Value = "C"
$ws.Cells.Item(8, 3).Value = "Miami Heat"
$ws.Cells.Item(9, 1).Value = "Julius Randle"
$ws.Cells.Item(9, 2).Value = "PF,C"
$ws.Cells.Item(9, 3).Value = "Minnesota Timberwolves"
$ws.Cells.Item(10, 1).Value = "Anthony Davis"
$ws.Cells.Item(10, 2).Value = "PF,C"
$ws.Cells.Item(10, 3).Value = "Los Angeles Lakers"
$ws.Cells.Item(11, 1).Value = "Herbert Jones"
$ws.Cells.Item(11, 2).Value = "SF,PF"
$ws.Cells.Item(11, 3).Value = "New Orleans Pelicans"
$ws.Cells.Item(12, 1).Value = "Isaiah Hartenstein"
$ws.Cells.Item(12, 2).Value = "C"
$ws.Cells.Item(12, 3).Value = "Oklahoma City Thunder"
$ws.Cells.Item(13, 1).Value = "Damian Lillard"
$ws.Cells.Item(13, 2).Value = "PG"
$ws.Cells.Item(13, 3).Value = "Milwaukee Bucks"
$ws.Cells.Item(14, 1).Value = "Bilal Coulibaly"
$ws.Cells.Item(14, 2).Value = "SG,SF"
$ws.Cells.Item(14, 3).Value = "Washington Wizards"
$ws.Cells.Item(15, 1).Value = "Derrick White"
$ws.Cells.Item(15, 2).Value = "PG,SG"
$ws.Cells.Item(15, 3).Value = "Boston Celtics"
$ws.Cells.Item(16, 1).Value = "Brandon Ingram"
$ws.Cells.Item(16, 2).Value = "SG,SF,PF"
$ws.Cells.Item(16, 3).Value = "New Orleans Pelicans"
$ws.Cells.Item(17, 1).Value = "LaMelo Ball"
$ws.Cells.Item(17, 2).Value = "PG,SG"
$ws.Cells.Item(17, 3).Value = "Charlotte Hornets"
